$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new "2022-Q4" row right after the
#    header, shifting the existing quarters down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# A9 doesn't exist yet - clone A2's formatting onto it before anything else
# moves around, so the new bottom row keeps the bold/border index style.
$total.Range("A2").Copy()
$total.Range("A9").PasteSpecial(-4122)

# Shift the quarter/count/value columns down one row (bottom-up so we never
# clobber a row before it has been read).
for ($r = 8; $r -ge 2; $r--) {
    $dst = $r + 1
    $total.Range("B$dst").Value = $total.Range("B$r").Value2
    $total.Range("C$dst").Value = $total.Range("C$r").Value2
    $total.Range("D$dst").Value = $total.Range("D$r").Value2
}

# New first data row: 2022-Q4
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 0.02

# Re-number the index column (0-based row offset) for every data row.
for ($r = 2; $r -le 9; $r++) {
    $total.Range("A$r").Value = $r - 2
}

# ---------------------------------------------------------------------------
# 2) Add a new "2022-Q4" fund-holdings sheet right after "总计", built from
#    a copy of "2022-Q3" (same layout/styles), then overwrite its data.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $total)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Each row: code, name, size, position%, ratio, marketValue, rank
$rows = @(
    @("004532", "民生加银中证港股通高股息精选指数A", "0.14", "92.86", "6.23", "0.0087", 3),
    @("004533", "民生加银中证港股通高股息精选指数C", "0.09", "92.86", "6.23", "0.0056", 3),
    @("006658", "财通中证香港红利等权投资指数A",       "0.14", "89.84", "3.17", "0.0044", 4),
    @("501307", "银河中证沪港深高股息指数（LOF）A",     "0.16", "93.15", "1.63", "0.0026", 3),
    @("006659", "财通中证香港红利等权投资指数C",       "0.04", "89.84", "3.17", "0.0013", 4),
    @("005770", "信澳中证沪港深高股息精选指数",         "0.13", "23.47", "0.69", "0.0009", 3),
    @("501308", "银河中证沪港深高股息指数（LOF）C",     "0.01", "93.15", "1.63", "0.0002", 3)
)

$r = 2
foreach ($row in $rows) {
    $q4.Range("B$r").Value = "'" + $row[0]
    $q4.Range("C$r").Value = $row[1]
    $q4.Range("D$r").Value = "'" + $row[2]
    $q4.Range("E$r").Value = "'" + $row[3]
    $q4.Range("F$r").Value = "'" + $row[4]
    $q4.Range("G$r").Value = "'" + $row[5]
    $q4.Range("H$r").Value = $row[6]
    $r = $r + 1
}

# Restore the original active tab (copying/renaming sheets along the way
# leaves the freshly-inserted sheet selected).
$total.Activate()
